$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.165.02'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = '1.787.38'
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').Value = '226.24'
$ws.Range('E5').Value = '  -0.52%  '
$ws.Range('D6').Value = '0.548'
$ws.Range('E6').Value = '  +0.57%  '
$ws.Range('D8').Value = '32.36'
$ws.Range('E8').Value = '  +0.19%  '
$ws.Range('E9').Value = '  -0.28%  '
$ws.Range('E10').Value = '  -0.09%  '
$ws.Range('E11').Value = '  +0.95%  '
$ws.Range('D12').Value = '2.044.99'
$ws.Range('E12').Value = '  -0.09%  '
$ws.Range('D13').Value = '11.05'
$ws.Range('E13').Value = '  -2.35%  '
$ws.Range('D14').Value = '1.790.36'
$ws.Range('E14').Value = '  +0.14%  '
$ws.Range('E15').Value = '  +0.40%  '
$ws.Range('D16').Value = '34.169.52'
$ws.Range('E16').Value = '  +0.34%  '
$ws.Range('E17').Value = '  +0.15%  '
$ws.Range('D18').Value = '67.84'
$ws.Range('E18').Value = '  -0.29%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').Value = '245.75'
$ws.Range('E19').Value = '  +1.05%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0797'
$ws.Range('E20').Value = '  +2.04%  '
$ws.Range('D21').Value = '11.04'
$ws.Range('E21').Value = '  +1.49%  '
$ws.Range('E22').Value = '  +0.16%  '
$ws.Range('D23').Value = '4.14'
$ws.Range('E23').Value = '  +0.89%  '
$ws.Range('D24').Value = '2.04'
$ws.Range('E24').Value = '  -0.60%  '
$ws.Range('D25').Value = '161.91'
$ws.Range('E25').Value = '  +0.44%  '
$ws.Range('E26').Value = '  -0.53%  '
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('E28').Value = '  +1.08%  '
$ws.Range('E29').Value = '  +0.39%  '
$ws.Range('E30').Value = '  -0.62%  '
$ws.Range('E31').Value = '  -0.27%  '
$ws.Range('E32').Value = '  +2.91%  '
$ws.Range('D33').Value = '3.73'
$ws.Range('E33').Value = '  +3.10%  '
$ws.Range('E34').Value = '  -1.60%  '
$ws.Range('D35').Value = '1.442.23'
$ws.Range('E35').Value = '  +2.32%  '
$ws.Range('D36').Value = '2.60'
$ws.Range('E36').Value = '  +10.10%  '
$ws.Range('D37').Value = '0.656'
$ws.Range('E37').Value = '  +1.04%  '
$ws.Range('E38').Value = '  +0.54%  '
$ws.Range('E39').Value = '  +0.63%  '
$ws.Range('D40').Value = '81.84'
$ws.Range('E40').Value = '  +1.84%  '
$ws.Range('E41').Value = '  +1.62%  '
$ws.Range('D42').Value = '13.85'
$ws.Range('E42').Value = '  +3.90%  '
$ws.Range('E43').Value = '  +1.02%  '
$ws.Range('D44').Value = '0.920'
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').Value = '0.0516'
$ws.Range('E45').Value = '  +1.53%  '
$ws.Range('E46').Value = '  +1.01%  '
$ws.Range('E47').Value = '  +0.79%  '
$ws.Range('D48').Value = '1.943.95'
$ws.Range('E48').Value = '  -0.17%  '
$ws.Range('D49').Value = '104.93'
$ws.Range('E49').Value = '  -1.79%  '
$ws.Range('E50').Value = '  +0.20%  '
$ws.Range('D51').Value = '0.0₆0130'
$ws.Range('E51').Value = '  -6.50%  '
